$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new time-tracking entry on row 21
$ws.Range("B21").Value = 45212
$ws.Range("C21").Value = 0.4548611111111111
$ws.Range("D21").Value = 0.56597222222222221
$ws.Range("D21").NumberFormat = $ws.Range("D20").NumberFormat
$ws.Range("G21").Value = "Got live editing working instead of doing what I said I needed to do next! Also did some theme tweaking"

# Match row height used for the wrapped text in row 21 (45pt, same pattern as other entry rows)
$ws.Rows.Item(21).RowHeight = 45

# Update the view to reflect where the author was scrolled to / selected when saving
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("G22").Select()
